$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric ("1.00", "6.06", etc.) but must remain
# stored as text (matching the source data, which is all inlineStr). Mark them
# as Text before writing, then clear the number format again so no extra
# cell style sticks around afterward.
$textCells = @('D4', 'D5', 'D6', 'D14', 'D18', 'D19', 'D21', 'D22', 'D24', 'D26', 'D27', 'D30', 'D31', 'D32', 'D33', 'D34', 'D36', 'D40', 'D41', 'D42', 'D44', 'D46', 'D47', 'D48', 'D49', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '71.491.53'
$ws.Range('E2').Value = '  +1.54%  '
$ws.Range('D3').Value = '3.823.83'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '700.69'
$ws.Range('E5').Value = '  +4.92%  '
$ws.Range('D6').Value = '174.12'
$ws.Range('E6').Value = '  +3.39%  '
$ws.Range('D7').Value = '3.819.68'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('E11').Value = '  +5.50%  '
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('E13').Value = '  +6.17%  '
$ws.Range('D14').Value = '36.63'
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('D15').Value = '4.471.15'
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').Value = '3.823.78'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').Value = '71.499.66'
$ws.Range('E17').Value = '  +1.64%  '
$ws.Range('D18').Value = '17.78'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').Value = '7.24'
$ws.Range('E19').Value = '  +1.19%  '
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('D21').Value = '11.16'
$ws.Range('E21').Value = '  +2.98%  '
$ws.Range('D22').Value = '487.84'
$ws.Range('E22').Value = '  +2.89%  '
$ws.Range('E23').Value = '  +0.73%  '
$ws.Range('D24').Value = '84.76'
$ws.Range('E24').Value = '  +2.45%  '
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('D26').Value = '12.37'
$ws.Range('E26').Value = '  +1.22%  '
$ws.Range('D27').Value = '10.57'
$ws.Range('E27').Value = '  +2.32%  '
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('D29').Value = '3.971.78'
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '3.15'
$ws.Range('E31').Value = '  +10.13%  '
$ws.Range('D32').Value = '2.31'
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('D33').Value = '7.60'
$ws.Range('E33').Value = '  +2.88%  '
$ws.Range('D34').Value = '29.73'
$ws.Range('E34').Value = '  +0.23%  '
$ws.Range('E35').Value = '  +2.03%  '
$ws.Range('D36').Value = '9.33'
$ws.Range('E36').Value = '  +1.95%  '
$ws.Range('D37').Value = '3.775.75'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('E39').Value = '  +2.00%  '
$ws.Range('D40').Value = '2.41'
$ws.Range('E40').Value = '  +15.17%  '
$ws.Range('D41').Value = '3.43'
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').Value = '6.06'
$ws.Range('E42').Value = '  +1.61%  '
$ws.Range('E43').Value = '  +3.13%  '
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D46').Value = '0.000312'
$ws.Range('E46').Value = '  +9.05%  '
$ws.Range('D47').Value = '163.23'
$ws.Range('E47').Value = '  +3.89%  '
$ws.Range('D48').Value = '44.82'
$ws.Range('E48').Value = '  -2.32%  '
$ws.Range('D49').Value = '48.73'
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('E50').Value = '  +0.71%  '
$ws.Range('D51').Value = '8.70'
$ws.Range('E51').Value = '  +2.40%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
